# Auto update Excel log 2026-02-04 14:10:18
# Appends newly captured sensor readings to the PIR, Humidity and Temperature
# sensor logs.  Column A (Date) and, for Humidity, column E (a "NN.N%" value)
# look like dates/percentages to Excel's automatic type inference, so those
# cells are briefly forced to Text format while the value is written and then
# reset to the workbook's default ("Normal") style, to match the existing
# plain-text cells already in the sheet.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($ws, $rows)

    foreach ($row in $rows) {
        $r = $row[0]

        # Date column - force text so "2026-02-04" isn't parsed as a date serial.
        $dateCell = $ws.Cells.Item($r, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $row[1]
        $dateCell.Style = "Normal"

        $ws.Cells.Item($r, 2).Value = $row[2]
        $ws.Cells.Item($r, 3).Value = $row[3]
        $ws.Cells.Item($r, 4).Value = $row[4]

        # Value column - force text so values like "78.6%" aren't parsed as numbers.
        $valueCell = $ws.Cells.Item($r, 5)
        $valueCell.NumberFormat = "@"
        $valueCell.Value = $row[5]
        $valueCell.Style = "Normal"

        $ws.Cells.Item($r, 6).Value = $row[6]
    }
}

$pirRows = @(
    @("97", "2026-02-04", "14:09:13", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("98", "2026-02-04", "14:09:14", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("99", "2026-02-04", "14:09:15", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("100", "2026-02-04", "14:09:16", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("101", "2026-02-04", "14:09:17", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("102", "2026-02-04", "14:09:22", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("103", "2026-02-04", "14:09:23", "14:00", "Bathroom", "Motion Detected", "Active"),
    @("104", "2026-02-04", "14:09:30", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("105", "2026-02-04", "14:09:35", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("106", "2026-02-04", "14:09:40", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("107", "2026-02-04", "14:09:45", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("108", "2026-02-04", "14:09:51", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("109", "2026-02-04", "14:09:56", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("110", "2026-02-04", "14:10:01", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("111", "2026-02-04", "14:10:06", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("112", "2026-02-04", "14:10:11", "14:00", "Bathroom", "No Motion", "Inactive")
)

$humidityRows = @(
    @("75", "2026-02-04", "14:09:13", "14:00", "Bathroom", "78.6%", "Active"),
    @("76", "2026-02-04", "14:09:14", "14:00", "Bathroom", "77.0%", "Active"),
    @("77", "2026-02-04", "14:09:14", "14:00", "Bathroom", "78.2%", "Active"),
    @("78", "2026-02-04", "14:09:15", "14:00", "Bathroom", "76.7%", "Active"),
    @("79", "2026-02-04", "14:09:16", "14:00", "Bathroom", "77.8%", "Active"),
    @("80", "2026-02-04", "14:09:19", "14:00", "Bathroom", "77.3%", "Active"),
    @("81", "2026-02-04", "14:09:30", "14:00", "Bathroom", "76.8%", "Active"),
    @("82", "2026-02-04", "14:09:35", "14:00", "Bathroom", "77.6%", "Active"),
    @("83", "2026-02-04", "14:09:40", "14:00", "Bathroom", "76.6%", "Active"),
    @("84", "2026-02-04", "14:09:45", "14:00", "Bathroom", "77.6%", "Active"),
    @("85", "2026-02-04", "14:09:55", "14:00", "Bathroom", "77.8%", "Active"),
    @("86", "2026-02-04", "14:10:00", "14:00", "Bathroom", "76.9%", "Active"),
    @("87", "2026-02-04", "14:10:05", "14:00", "Bathroom", "77.8%", "Active"),
    @("88", "2026-02-04", "14:10:10", "14:00", "Bathroom", "76.9%", "Active")
)

$temperatureRows = @(
    @("75", "2026-02-04", "14:09:13", "14:00", "Bathroom", "24.8C", "Active"),
    @("76", "2026-02-04", "14:09:14", "14:00", "Bathroom", "24.9C", "Active"),
    @("77", "2026-02-04", "14:09:15", "14:00", "Bathroom", "24.8C", "Active"),
    @("78", "2026-02-04", "14:09:16", "14:00", "Bathroom", "24.8C", "Active"),
    @("79", "2026-02-04", "14:09:17", "14:00", "Bathroom", "24.8C", "Active"),
    @("80", "2026-02-04", "14:09:20", "14:00", "Bathroom", "24.9C", "Active"),
    @("81", "2026-02-04", "14:09:30", "14:00", "Bathroom", "24.8C", "Active"),
    @("82", "2026-02-04", "14:09:35", "14:00", "Bathroom", "24.8C", "Active"),
    @("83", "2026-02-04", "14:09:40", "14:00", "Bathroom", "24.8C", "Active"),
    @("84", "2026-02-04", "14:09:45", "14:00", "Bathroom", "24.8C", "Active"),
    @("85", "2026-02-04", "14:09:55", "14:00", "Bathroom", "24.8C", "Active"),
    @("86", "2026-02-04", "14:10:00", "14:00", "Bathroom", "24.8C", "Active"),
    @("87", "2026-02-04", "14:10:05", "14:00", "Bathroom", "24.8C", "Active"),
    @("88", "2026-02-04", "14:10:10", "14:00", "Bathroom", "24.8C", "Active")
)

$wsPir = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPir $pirRows

$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity $humidityRows

$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature $temperatureRows
